$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row containing account "005040864" / "ANDRE" / 10142.82
# This is row 7 (row 1 is the header: Conta, Nome, Saldo)
$ws.Rows.Item(7).Delete()
